$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new data rows at row 200, pushing the existing rows
# (old 200..267) down to 202..269.
$ws.Rows.Item(200).Insert()
$ws.Rows.Item(200).Insert()

# Populate the first new row (200)
$ws.Range("A200").Value = 11
$ws.Range("B200").Value = 'Vega Monumental Concepción'
$ws.Range("C200").Value = 'Bíobío'
$ws.Range("D200").Value = 44642
$ws.Range("E200").Value = 8
$ws.Range("F200").Value = 100112002
$ws.Range("G200").Value = 'Pimiento'
$ws.Range("H200").Value = 'Cuatro cascos verde'
$ws.Range("I200").Value = 'Primera'
$ws.Range("J200").Value = 180
$ws.Range("K200").Value = 15000
$ws.Range("L200").Value = 16000
$ws.Range("M200").Value = 15444
$ws.Range("N200").Value = '$/caja 18 kilos'
$ws.Range("O200").Value = 'Provincia de Limarí'
$ws.Range("P200").Value = 858
$ws.Range("Q200").Value = 18
$ws.Range("R200").Value = 'Hortaliza'

# Populate the second new row (201)
$ws.Range("A201").Value = 11
$ws.Range("B201").Value = 'Vega Monumental Concepción'
$ws.Range("C201").Value = 'Bíobío'
$ws.Range("D201").Value = 44642
$ws.Range("E201").Value = 8
$ws.Range("F201").Value = 100112002
$ws.Range("G201").Value = 'Pimiento'
$ws.Range("H201").Value = 'Morrón rojo'
$ws.Range("I201").Value = 'Primera'
$ws.Range("J201").Value = 180
$ws.Range("K201").Value = 15000
$ws.Range("L201").Value = 16000
$ws.Range("M201").Value = 15444
$ws.Range("N201").Value = '$/caja 18 kilos'
$ws.Range("O201").Value = 'Provincia de Limarí'
$ws.Range("P201").Value = 858
$ws.Range("Q201").Value = 18
$ws.Range("R201").Value = 'Hortaliza'
